$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the bilinear interpolation formula: the second squared term should be
#    ADDED, not subtracted, so the tests actually exercise the interpolation
#    properly. Re-seed both shared-formula "roots" (B2, and C2:L12 / B3:B12)
#    and let Excel propagate + recalc the relative references across the
#    range, which reconstructs the same shared-formula grouping as before.
$ws.Range("B2").Formula = "=EXP(-(B`$1^2)+(`$A2^2))"
$ws.Range("C2:L12").Formula = "=EXP(-(C`$1^2)+(`$A2^2))"
$ws.Range("B3:B12").Formula = "=EXP(-(B`$1^2)+(`$A3^2))"

# 2. Apply a fixed-point numeric display format (9 decimal digits) to the
#    interpolated results grid.
$ws.Range("B2:L12").NumberFormat = "0.000000000"

# 3. Resize the columns to fit their new (wider/longer) displayed content.
$ws.Columns("A").ColumnWidth = 8.0
$ws.Columns("B").ColumnWidth = 8.666666666666666
$ws.Columns("C:L").ColumnWidth = 8.5

# 4. Update the active selection to cover the results grid.
$ws.Range("B2:L12").Select() | Out-Null
